# Update Efnb3-Rhbdl2 LR-pair data with new TPM values, adding Resolving-Mac rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Rhbdl2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1533166666666667
$ws.Range("H2").Value = 0.45995
$ws.Range("I2").Value = 0.1101680340964018
$ws.Range("J2").Value = 0.1101680340964017
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.082435
$ws.Range("N2").Value = 6.247305
$ws.Range("O2").Value = 0.9920259111440977
$ws.Range("P2").Value = 0.9920259111440977
$ws.Range("Q2").Value = 0.31927199275
$ws.Range("R2").Value = 2.87344793475
$ws.Range("S2").Value = 0.109289544403437
$ws.Range("T2").Value = 0.109289544403437

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Rhbdl2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1533166666666667
$ws.Range("H3").Value = 0.45995
$ws.Range("I3").Value = 0.1101680340964018
$ws.Range("J3").Value = 0.1101680340964017
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.016739
$ws.Range("N3").Value = 0.050217
$ws.Range("O3").Value = 0.00797408885590237
$ws.Range("P3").Value = 0.007974088855902369
$ws.Range("Q3").Value = 0.002566367683333334
$ws.Range("R3").Value = 0.02309730915
$ws.Range("S3").Value = 0.0008784896929647897
$ws.Range("T3").Value = 0.0008784896929647894

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Rhbdl2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9559960000000002
$ws.Range("H4").Value = 2.867988
$ws.Range("I4").Value = 0.6869455370628789
$ws.Range("J4").Value = 0.6869455370628788
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.082435
$ws.Range("N4").Value = 6.247305
$ws.Range("O4").Value = 0.9920259111440977
$ws.Range("P4").Value = 0.9920259111440977
$ws.Range("Q4").Value = 1.99079953026
$ws.Range("R4").Value = 17.91719577234
$ws.Range("S4").Value = 0.681467772311174
$ws.Range("T4").Value = 0.6814677723111738

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb3"
$ws.Range("C5").Value = "Rhbdl2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9559960000000002
$ws.Range("H5").Value = 2.867988
$ws.Range("I5").Value = 0.6869455370628789
$ws.Range("J5").Value = 0.6869455370628788
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.016739
$ws.Range("N5").Value = 0.050217
$ws.Range("O5").Value = 0.00797408885590237
$ws.Range("P5").Value = 0.007974088855902369
$ws.Range("Q5").Value = 0.016002417044
$ws.Range("R5").Value = 0.144021753396
$ws.Range("S5").Value = 0.005477764751704971
$ws.Range("T5").Value = 0.00547776475170497

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Efnb3"
$ws.Range("C6").Value = "Rhbdl2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1265133333333333
$ws.Range("H6").Value = 0.37954
$ws.Range("I6").Value = 0.09090808927263468
$ws.Range("J6").Value = 0.09090808927263468
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.082435
$ws.Range("N6").Value = 6.247305
$ws.Range("O6").Value = 0.9920259111440977
$ws.Range("P6").Value = 0.9920259111440977
$ws.Range("Q6").Value = 0.2634557933
$ws.Range("R6").Value = 2.3711021397
$ws.Range("S6").Value = 0.09018318009105439
$ws.Range("T6").Value = 0.09018318009105439

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Efnb3"
$ws.Range("C7").Value = "Rhbdl2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1265133333333333
$ws.Range("H7").Value = 0.37954
$ws.Range("I7").Value = 0.09090808927263468
$ws.Range("J7").Value = 0.09090808927263468
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.016739
$ws.Range("N7").Value = 0.050217
$ws.Range("O7").Value = 0.00797408885590237
$ws.Range("P7").Value = 0.007974088855902369
$ws.Range("Q7").Value = 0.002117706686666667
$ws.Range("R7").Value = 0.01905936018
$ws.Range("S7").Value = 0.000724909181580294
$ws.Range("T7").Value = 0.0007249091815802938

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Efnb3"
$ws.Range("C8").Value = "Rhbdl2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.155836
$ws.Range("H8").Value = 0.467508
$ws.Range("I8").Value = 0.1119783395680848
$ws.Range("J8").Value = 0.1119783395680848
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.082435
$ws.Range("N8").Value = 6.247305
$ws.Range("O8").Value = 0.9920259111440977
$ws.Range("P8").Value = 0.9920259111440977
$ws.Range("Q8").Value = 0.32451834066
$ws.Range("R8").Value = 2.92066506594
$ws.Range("S8").Value = 0.1110854143384325
$ws.Range("T8").Value = 0.1110854143384324

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Efnb3"
$ws.Range("C9").Value = "Rhbdl2"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.155836
$ws.Range("H9").Value = 0.467508
$ws.Range("I9").Value = 0.1119783395680848
$ws.Range("J9").Value = 0.1119783395680848
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.016739
$ws.Range("N9").Value = 0.050217
$ws.Range("O9").Value = 0.00797408885590237
$ws.Range("P9").Value = 0.007974088855902369
$ws.Range("Q9").Value = 0.002608538804
$ws.Range("R9").Value = 0.023476849236
$ws.Range("S9").Value = 0.0008929252296523162
$ws.Range("T9").Value = 0.000892925229652316
